# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-06 (serial 45175) to 2023-09-08 (serial 45177).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 111 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
